$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E header: "Matières enseignés" (adds a new shared string,
# extends the used range to A1:E1)
$ws.Range("E1").Value = "Matières enseignés"

# Column widths for C, D and E (values chosen so the engine's
# character->pixel rounding lands as close as possible to the authored
# widths of 27.5703125 / 15.7109375 / 31.7109375)
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 30.833333333333332

# Leave the selection on the new column, matching the saved session state
$ws.Range("E6").Select() | Out-Null
